$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in Home_Score (J) and Away_Score (K) for the fixtures played on Jun 19-20, 2024
$ws.Range("J15").Value = 2
$ws.Range("K15").Value = 0

$ws.Range("J16").Value = 1
$ws.Range("K16").Value = 1

$ws.Range("J17").Value = 1
$ws.Range("K17").Value = 1

$ws.Range("J18").Value = 1
$ws.Range("K18").Value = 1

$ws.Range("J19").Value = 1
$ws.Range("K19").Value = 0
